# Update "想去人数" (F column) figures that changed between the previous
# and current data pulls, on both the "展览" sheet and the aggregated
# "全部类型" sheet.

$wb = $excel.ActiveWorkbook

# Sheet "展览" — row => new F value
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value  = 68
$ws1.Range("F7").Value  = 2663
$ws1.Range("F9").Value  = 250
$ws1.Range("F10").Value = 106
$ws1.Range("F11").Value = 9818
$ws1.Range("F12").Value = 70
$ws1.Range("F13").Value = 252
$ws1.Range("F14").Value = 3
$ws1.Range("F15").Value = 603
$ws1.Range("F16").Value = 11701
$ws1.Range("F17").Value = 12021
$ws1.Range("F19").Value = 86
$ws1.Range("F21").Value = 22

# Sheet "全部类型" — same events, shifted down one row because it also
# includes the single "演出" entry.
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value  = 68
$ws4.Range("F7").Value  = 2663
$ws4.Range("F10").Value = 250
$ws4.Range("F11").Value = 106
$ws4.Range("F12").Value = 9818
$ws4.Range("F13").Value = 70
$ws4.Range("F14").Value = 252
$ws4.Range("F15").Value = 3
$ws4.Range("F16").Value = 603
$ws4.Range("F17").Value = 11701
$ws4.Range("F18").Value = 12021
$ws4.Range("F20").Value = 87
$ws4.Range("F22").Value = 22
